# "implement and validate sql template"
#
# On the DataTypes sheet, the T-SQL column used the fixed-length placeholder
# "(0)" for variable/blob-ish text types. Update it to use SQL Server's
# MAX-length syntax instead:
#   - blob    -> T-SQL VARBINARY(0)  => VARBINARY(MAX)
#   - char    -> T-SQL CHAR(0)       => CHAR(MAX)
#   - varchar -> T-SQL VARCHAR(0)    => VARCHAR(MAX)
#
# (The MySQL/SQLite columns for char & varchar keep their existing values —
# only the T-SQL column changes.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataTypes")

$ws.Range("E3").Value = "VARBINARY(MAX)"   # blob  / T-SQL
$ws.Range("E4").Value = "CHAR(MAX)"        # char  / T-SQL
$ws.Range("E6").Value = "VARCHAR(MAX)"     # varchar / T-SQL

# Move the active selection to match the saved workbook state.
$ws.Range("E7").Select()
